$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.351.88"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.300.42"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "588.85"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").Value = "179.80"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").Value = "0.642"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "3.291.37"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "6.83"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "0.401"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "3.869.53"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "66.323.32"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "26.59"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "3.307.01"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "426.70"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").Value = "13.01"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").Value = "7.30"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "71.26"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").Value = "5.67"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "0.511"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.205"
$ws.Range("E27").Value = "  +5.51%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0000114"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "9.36"
$ws.Range("E29").Value = "  +4.77%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.92"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "22.27"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "5.16"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "6.55"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "1.18"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "158.96"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.43"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.851.50"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "26.20"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "4.32"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.746"
$ws.Range("E43").Value = "  -4.55%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "39.62"
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "5.88"
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.31"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0639"
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "313.73"
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "22.78"
$ws.Range("E49").Value = "  -3.00%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0270"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.102"
$ws.Range("E51").Value = "  -0.06%  "
